# Commit: "pulled date range creation to it's own function"
#
# Net effect observed in the target workbook:
#   * the existing "inflation" worksheet was renamed to "related" and now
#     holds a small lookup/expense table (position_id..end_date)
#   * a brand-new "inflation" worksheet was (re)created right after it,
#     carrying forward exactly what the old "inflation" sheet used to
#     contain (the EDATE/compounding schedule)
#   * positions!F6 (salary_annual for position 5) changed 95000 -> 100000
#
# We reproduce this with the Excel object model by literally copying the
# current "inflation" sheet (so the new sheet keeps its formulas/formats),
# then renaming the original into "related" and repurposing its cells for
# the new table, and the copy back into "inflation".

$wb = $excel.ActiveWorkbook

$positions = $wb.Worksheets.Item("positions")
$inflation = $wb.Worksheets.Item("inflation")

# 1) Duplicate "inflation" immediately after itself - the duplicate will
#    become the "new" inflation sheet, the original becomes "related".
#    Fetch it back by position (rather than assuming Excel's generated
#    "inflation (2)" name) so this doesn't depend on naming/locale quirks.
$inflationIndex = $inflation.Index
$inflation.Copy($null, $inflation)
$inflationCopy = $wb.Worksheets.Item($inflationIndex + 1)

$inflation.Name = "related"
$inflationCopy.Name = "inflation"

$related = $wb.Worksheets.Item("related")

# 2) Wipe the old inflation-schedule content/format out of "related" - it's
#    going to hold a different table now.
$related.Cells.Clear()

# 3) Bring over the plain "data" style (font, no special number format)
#    used throughout the rest of the workbook for the bulk of the table...
$positions.Range("A1").Copy()
$related.Range("A1:J4").PasteSpecial(-4122)

# ... and the date style (m/d/yyyy) for the date columns.
$positions.Range("I3").Copy()
$related.Range("I2:I4").PasteSpecial(-4122)
$related.Range("J4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# end_date is entered with a shorter m/d/yy format for the first two rows.
$related.Range("J2:J3").NumberFormat = "m/d/yy"

# 4) Header row.
$headers = @("position_id", "position_title", "department", "employee_id", "employee_name", "item", "expense_type", "amount_annual", "start_date", "end_date")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $related.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# 5) Table rows.
$related.Range("A2").Value = 1
$related.Range("B2").Value = "CEO"
$related.Range("C2").Value = "d1"
$related.Range("D2").Value = "e10001"
$related.Range("E2").Value = "Alice"
$related.Range("F2").Value = "Stock compensation"
$related.Range("G2").Value = "salary"
$related.Range("H2").Value = 10000
$related.Range("I2").Value = 44562
$related.Range("J2").Value = 45291

$related.Range("A3").Value = 2
$related.Range("B3").Value = "CFO"
$related.Range("C3").Value = "d2"
$related.Range("D3").Value = "e10002"
$related.Range("E3").Value = "Bob"
$related.Range("F3").Value = "Home office allowance"
$related.Range("G3").Value = "office expense"
$related.Range("H3").Value = 1200
$related.Range("I3").Value = 44576
$related.Range("J3").Value = 44972

$related.Range("A4").Value = 3
$related.Range("B4").Value = "COO"
$related.Range("C4").Value = "d3"
$related.Range("D4").Value = "e10003"
$related.Range("E4").Value = "Charlie"
$related.Range("F4").Value = "Travel"
$related.Range("G4").Value = "travel"
$related.Range("H4").Value = 5000
$related.Range("I4").Value = 44561
$related.Range("J4").Value = 44654

# 6) The unrelated salary bump for position 5 on the "positions" sheet.
$positions.Range("F6").Value = 100000
